# Auto-generated script to apply cryptos.xlsx (coinranking) update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.359.51'
$ws.Range("E2").Value = '  -1.37%  '
$ws.Range("D3").Value = '1.593.20'
$ws.Range("E3").Value = '  -0.43%  '
$ws.Range("E4").Value = '  -0.47%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '210.43'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -0.60%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.505'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -1.37%  '
$ws.Range("E8").Value = '  -1.10%  '
$ws.Range("E9").Value = '  -0.43%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.65'
$ws.Range("D10").ClearFormats()
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0845'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  -0.38%  '
$ws.Range("D12").Value = '1.814.64'
$ws.Range("E12").Value = '  -0.50%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.07'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  +0.21%  '
$ws.Range("D14").Value = '1.585.13'
$ws.Range("E14").Value = '  -1.01%  '
$ws.Range("E15").Value = '  -1.16%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '64.70'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  -0.57%  '
$ws.Range("D17").Value = '26.362.89'
$ws.Range("E17").Value = '  -1.22%  '
$ws.Range("D18").Value = '0.0₃0730'
$ws.Range("E18").Value = '  -1.50%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.49'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  +4.75%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '212.26'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  +1.57%  '
$ws.Range("E21").Value = '  -0.48%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.29'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -0.23%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '2.20'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -1.73%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '8.94'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  -0.89%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '145.19'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +1.07%  '
$ws.Range("E26").Value = '  -0.42%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.07'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -0.83%  '
$ws.Range("E28").Value = '  -1.02%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.30'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -0.41%  '
$ws.Range("E30").Value = '  -0.54%  '
$ws.Range("E31").Value = '  -0.78%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.24'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -0.76%  '
$ws.Range("E33").Value = '  +0.23%  '
$ws.Range("D34").Value = '1.301.86'
$ws.Range("E34").Value = '  +1.82%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.618'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  +4.51%  '
$ws.Range("E36").Value = '  -1.76%  '
$ws.Range("E37").Value = '  -1.08%  '
$ws.Range("E38").Value = '  -0.55%  '
$ws.Range("E39").Value = '  -13.06%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.817'
$ws.Range("D40").ClearFormats()
$ws.Range("E41").Value = '  -0.52%  '
$ws.Range("E42").Value = '  +2.74%  '
$ws.Range("B43").Value = 'Aave'
$ws.Range("C43").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '62.88'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  +0.26%  '
$ws.Range("B44").Value = 'MXToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.14'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -2.35%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.763'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -1.94%  '
$ws.Range("D46").Value = '1.728.49'
$ws.Range("E46").Value = '  -0.35%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '88.59'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -2.14%  '
$ws.Range("E48").Value = '  -3.75%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0994'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -2.70%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0506'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -1.21%  '
$ws.Range("B51").Value = 'EnergySwap'
$ws.Range("C51").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.45'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -1.56%  '
